$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Worksheet 1"

# --- Row 2 (new data row) ---
$ws.Range("A2").Value = "stuff"
$ws.Range("B2").Value = "things"
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = "have"
$ws.Range("E2").Value = "you"
$ws.Range("F2").Value = "to"
$ws.Range("G2").Value = "say"
$ws.Range("H2").Value = $true
$ws.Range("I2").Value = $false
$ws.Range("J2").Value = "for"
$ws.Range("I1").Value = " fatal"
$ws.Range("K2").Value = "your"
$ws.Range("L2").Value = "self"
$ws.Range("M2").Value = "ssp"
$ws.Range("N2").Value = "adult"
$ws.Range("T2").Value = "plk"

# --- Row 3 (new data row) ---
$ws.Range("A3").Value = "stuff"
$ws.Range("B3").Value = "things"
$ws.Range("D3").Value = "have"
$ws.Range("E3").Value = "you"
$ws.Range("F3").Value = "to"
$ws.Range("G3").Value = "say"
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = $false
$ws.Range("J3").Value = "for"
$ws.Range("K3").Value = "your"
$ws.Range("L3").Value = "self"
$ws.Range("N3").Value = "adult"
$ws.Range("T3").Value = "plk"

# --- Header row: shorten the diseaseDetected header (written last) ---
$ws.Range("H1").Value = " diseaseDetected"

# Update the selected cell to match the saved view state
$sel = $ws.Range("M2").Select()
